$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $savedStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $savedStyle
}

Set-TextValue 'D2' '68.134.18'
Set-TextValue 'E2' '  +0.52%  '
Set-TextValue 'D3' '3.798.17'
Set-TextValue 'E3' '  -0.16%  '
Set-TextValue 'E4' '  -0.18%  '
Set-TextValue 'D5' '601.62'
Set-TextValue 'E5' '  +0.75%  '
Set-TextValue 'D6' '165.13'
Set-TextValue 'E6' '  -1.50%  '
Set-TextValue 'E7' '  -0.09%  '
Set-TextValue 'E8' '  -0.56%  '
Set-TextValue 'E10' '  +0.40%  '
Set-TextValue 'E11' '  +3.15%  '
Set-TextValue 'E12' '  -2.12%  '
Set-TextValue 'D13' '35.88'
Set-TextValue 'E13' '  -0.37%  '
Set-TextValue 'D14' '4.434.55'
Set-TextValue 'E14' '  -0.19%  '
Set-TextValue 'D15' '3.781.32'
Set-TextValue 'E15' '  -0.54%  '
Set-TextValue 'D16' '68.129.48'
Set-TextValue 'E16' '  +0.52%  '
Set-TextValue 'E17' '  -1.21%  '
Set-TextValue 'E18' '  +2.37%  '
Set-TextValue 'D19' '7.08'
Set-TextValue 'E19' '  -0.21%  '
Set-TextValue 'D20' '461.65'
Set-TextValue 'E20' '  -0.01%  '
Set-TextValue 'D21' '9.71'
Set-TextValue 'E21' '  -2.46%  '
Set-TextValue 'E22' '  +0.05%  '
Set-TextValue 'E23' '  -4.50%  '
Set-TextValue 'D24' '83.03'
Set-TextValue 'E24' '  -0.68%  '
Set-TextValue 'D25' '12.01'
Set-TextValue 'E25' '  -0.65%  '
Set-TextValue 'E26' '  +0.06%  '
Set-TextValue 'B27' 'Dai'
Set-TextValue 'C27' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D27' '1.00'
Set-TextValue 'B28' 'RenderToken'
Set-TextValue 'C28' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D28' '10.02'
Set-TextValue 'D29' '3.945.66'
Set-TextValue 'E29' '  -0.23%  '
Set-TextValue 'E30' '  -0.18%  '
Set-TextValue 'E31' '  -4.63%  '
Set-TextValue 'E32' '  +0.77%  '
Set-TextValue 'D33' '29.35'
Set-TextValue 'E33' '  -1.01%  '
Set-TextValue 'D34' '1.00'
Set-TextValue 'D35' '9.04'
Set-TextValue 'E35' '  -0.62%  '
Set-TextValue 'D36' '0.0997'
Set-TextValue 'E36' '  -0.24%  '
Set-TextValue 'D37' '3.34'
Set-TextValue 'E37' '  -2.50%  '
Set-TextValue 'D39' '5.84'
Set-TextValue 'E39' '  +1.02%  '
Set-TextValue 'D40' '0.990'
Set-TextValue 'E40' '  -1.30%  '
Set-TextValue 'E41' '  +0.01%  '
Set-TextValue 'D43' '47.64'
Set-TextValue 'E43' '  -0.98%  '
Set-TextValue 'E44' '  +0.00%  '
Set-TextValue 'D45' '43.12'
Set-TextValue 'E45' '  -1.27%  '
Set-TextValue 'D46' '152.54'
Set-TextValue 'E46' '  +2.67%  '
Set-TextValue 'D47' '8.36'
Set-TextValue 'E47' '  +0.43%  '
Set-TextValue 'E48' '  +2.03%  '
Set-TextValue 'D49' '1.36'
Set-TextValue 'E49' '  +1.39%  '
Set-TextValue 'D50' '392.06'
Set-TextValue 'E50' '  -0.71%  '
Set-TextValue 'D51' '26.66'
Set-TextValue 'E51' '  -1.19%  '
